# Auto-generated edit script for cryptos.xlsx update
# Updates Price (D) and Volume(1h) (E) columns, and swaps the
# Chainlink/Polkadot rows (13 and 14) per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.026.57"
$ws.Range("E2").Value = "'  +2.54%  "
$ws.Range("D3").Value = "'1.820.63"
$ws.Range("E3").Value = "'  +3.00%  "
$ws.Range("E4").Value = "'  +0.94%  "
$ws.Range("D5").Value = "'314.75"
$ws.Range("E5").Value = "'  +3.09%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "'  +0.82%  "
$ws.Range("D7").Value = "'0.4313"
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("D8").Value = "'0.3700"
$ws.Range("E8").Value = "'  +1.77%  "
$ws.Range("D9").Value = "'0.07278"
$ws.Range("E9").Value = "'  +2.87%  "
$ws.Range("D10").Value = "'2.126.84"
$ws.Range("E10").Value = "'  +19.91%  "
$ws.Range("D11").Value = "'0.8704"
$ws.Range("E11").Value = "'  +2.19%  "
$ws.Range("E12").Value = "'  +5.29%  "
$ws.Range("B13").Value = "'Polkadot"
$ws.Range("C13").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.424"
$ws.Range("E13").Value = "'  +3.08%  "
$ws.Range("B14").Value = "'Chainlink"
$ws.Range("C14").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'6.653"
$ws.Range("E14").Value = "'  +3.31%  "
$ws.Range("D15").Value = "'0.06968"
$ws.Range("E15").Value = "'  +2.63%  "
$ws.Range("D16").Value = "'81.24"
$ws.Range("E16").Value = "'  +2.57%  "
$ws.Range("D17").Value = "'1.010"
$ws.Range("E17").Value = "'  +0.65%  "
$ws.Range("E18").Value = "'  +2.38%  "
$ws.Range("D19").Value = "'1.009"
$ws.Range("E19").Value = "'  +0.88%  "
$ws.Range("D20").Value = "'15.36"
$ws.Range("E20").Value = "'  +2.42%  "
$ws.Range("D21").Value = "'27.084.03"
$ws.Range("E21").Value = "'  +2.73%  "
$ws.Range("D22").Value = "'5.214"
$ws.Range("E22").Value = "'  +3.86%  "
$ws.Range("E23").Value = "'  -1.72%  "
$ws.Range("D24").Value = "'2.376.20"
$ws.Range("E24").Value = "'  +19.29%  "
$ws.Range("D25").Value = "'154.60"
$ws.Range("E25").Value = "'  +1.20%  "
$ws.Range("D26").Value = "'1.891"
$ws.Range("E26").Value = "'  +1.47%  "
$ws.Range("E27").Value = "'  +1.74%  "
$ws.Range("D28").Value = "'5.254"
$ws.Range("E28").Value = "'  +3.38%  "
$ws.Range("D29").Value = "'1.925"
$ws.Range("E29").Value = "'  +11.93%  "
$ws.Range("D30").Value = "'115.05"
$ws.Range("E30").Value = "'  +0.62%  "
$ws.Range("D31").Value = "'0.08991"
$ws.Range("E31").Value = "'  +0.65%  "
$ws.Range("D32").Value = "'1.177"
$ws.Range("E32").Value = "'  +5.50%  "
$ws.Range("D33").Value = "'0.7487"
$ws.Range("E33").Value = "'  +2.23%  "
$ws.Range("D34").Value = "'4.436"
$ws.Range("E34").Value = "'  +2.21%  "
$ws.Range("D35").Value = "'2.817"
$ws.Range("E35").Value = "'  +1.91%  "
$ws.Range("D36").Value = "'1.009"
$ws.Range("E36").Value = "'  +0.90%  "
$ws.Range("E37").Value = "'  +4.69%  "
$ws.Range("D38").Value = "'0.05253"
$ws.Range("E38").Value = "'  +2.21%  "
$ws.Range("D39").Value = "'0.01930"
$ws.Range("E39").Value = "'  +2.04%  "
$ws.Range("D40").Value = "'0.5123"
$ws.Range("E40").Value = "'  +4.03%  "
$ws.Range("D41").Value = "'2.756"
$ws.Range("E41").Value = "'  +9.60%  "
$ws.Range("D42").Value = "'0.1658"
$ws.Range("E42").Value = "'  +3.17%  "
$ws.Range("D43").Value = "'6.518"
$ws.Range("E43").Value = "'  +4.39%  "
$ws.Range("D44").Value = "'8.335"
$ws.Range("E44").Value = "'  +3.00%  "
$ws.Range("D45").Value = "'107.40"
$ws.Range("E45").Value = "'  +2.02%  "
$ws.Range("D46").Value = "'10.49"
$ws.Range("E46").Value = "'  +4.11%  "
$ws.Range("D47").Value = "'1.010"
$ws.Range("E47").Value = "'  +0.97%  "
$ws.Range("D48").Value = "'1.658"
$ws.Range("E48").Value = "'  +4.86%  "
$ws.Range("D49").Value = "'0.4578"
$ws.Range("E49").Value = "'  +1.96%  "
$ws.Range("D50").Value = "'0.06235"
$ws.Range("E50").Value = "'  +0.60%  "
$ws.Range("D51").Value = "'1.838"
$ws.Range("E51").Value = "'  +5.98%  "
